$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$oldText = $ws1.Range("A1").Value2
$newText = $oldText -replace [regex]::Escape("Bs = 12.63 = 51112.37 pesos"), "Bs = 12.5 = 50300.0 pesos"
$newText = $newText -replace [regex]::Escape("51112.37 pesos = 12.62 = 959.02 Bs"), "50300.0 pesos = 12.42 = 962.48 Bs"
$ws1.Range("A1").Value = $newText

# --- Update rate figures on "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 80
$ws2.Range("O10").Value = 4024
$ws2.Range("N12").Value = 4049
$ws2.Range("O12").Value = 77.477
